$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2094017094017094
$ws.Range("C2").Value = 0.5470085470085471
$ws.Range("J2").Value = 0.008547008547008548
$ws.Range("P2").Value = 0.1495726495726496
$ws.Range("S2").Value = 0.08547008547008547
$ws.Range("B3").Value = 0.007407407407407408
$ws.Range("C3").Value = 0.05185185185185185
$ws.Range("J3").Value = 0.02222222222222222
$ws.Range("P3").Value = 0.7185185185185186
$ws.Range("S3").Value = 0.2
$ws.Range("J4").Value = 0.09090909090909091
$ws.Range("P4").Value = 0.6060606060606061
$ws.Range("S4").Value = 0.303030303030303
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.06222222222222222
$ws.Range("D6").Value = 0.01333333333333333
$ws.Range("F6").Value = 0.07111111111111111
$ws.Range("J6").Value = 0.2355555555555555
$ws.Range("O6").Value = 0.02666666666666667
$ws.Range("Q6").Value = 0.1555555555555556
$ws.Range("R6").Value = 0.04444444444444445
$ws.Range("S6").Value = 0.3911111111111111
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.01515151515151515
$ws.Range("F7").Value = 0.08585858585858586
$ws.Range("J7").Value = 0.1464646464646465
$ws.Range("O7").Value = 0.02525252525252525
$ws.Range("Q7").Value = 0.1212121212121212
$ws.Range("R7").Value = 0.03535353535353535
$ws.Range("S7").Value = 0.4595959595959596
$ws.Range("B8").Value = 0.09024390243902439
$ws.Range("D8").Value = 0.01219512195121951
$ws.Range("E8").Value = 0.002439024390243902
$ws.Range("F8").Value = 0.07073170731707316
$ws.Range("J8").Value = 0.07560975609756097
$ws.Range("O8").Value = 0.02926829268292683
$ws.Range("Q8").Value = 0.1853658536585366
$ws.Range("R8").Value = 0.08780487804878048
$ws.Range("S8").Value = 0.4463414634146342
$ws.Range("B9").Value = 0.06842105263157895
$ws.Range("D9").Value = 0.02631578947368421
$ws.Range("F9").Value = 0.07368421052631578
$ws.Range("J9").Value = 0.1105263157894737
$ws.Range("O9").Value = 0.05263157894736842
$ws.Range("Q9").Value = 0.1526315789473684
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.4157894736842105
$ws.Range("B10").Value = 0.09176029962546817
$ws.Range("D10").Value = 0.01591760299625468
$ws.Range("F10").Value = 0.07303370786516854
$ws.Range("J10").Value = 0.1086142322097378
$ws.Range("O10").Value = 0.01685393258426966
$ws.Range("Q10").Value = 0.2134831460674157
$ws.Range("R10").Value = 0.08146067415730338
$ws.Range("S10").Value = 0.398876404494382
$ws.Range("F11").Value = 0.003311258278145695
$ws.Range("G11").Value = 0.1556291390728477
$ws.Range("J11").Value = 0.09271523178807947
$ws.Range("K11").Value = 0.2019867549668874
$ws.Range("L11").Value = 0.5364238410596026
$ws.Range("S11").Value = 0.009933774834437087
$ws.Range("G12").Value = 0.7455621301775148
$ws.Range("J12").Value = 0.1715976331360947
$ws.Range("K12").Value = 0.01775147928994083
$ws.Range("L12").Value = 0.04142011834319527
$ws.Range("S12").Value = 0.02366863905325444
$ws.Range("G13").Value = 0.6904761904761905
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.02380952380952381
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25
$ws.Range("F15").Value = 0.02358490566037736
$ws.Range("H15").Value = 0.1650943396226415
$ws.Range("I15").Value = 0.07547169811320754
$ws.Range("J15").Value = 0.3160377358490566
$ws.Range("K15").Value = 0.06132075471698113
$ws.Range("M15").Value = 0.01415094339622642
$ws.Range("O15").Value = 0.05660377358490566
$ws.Range("S15").Value = 0.2877358490566038
$ws.Range("F16").Value = 0.02
$ws.Range("H16").Value = 0.18
$ws.Range("I16").Value = 0.04666666666666667
$ws.Range("J16").Value = 0.4066666666666667
$ws.Range("K16").Value = 0.1066666666666667
$ws.Range("M16").Value = 0.01333333333333333
$ws.Range("N16").Value = 0.006666666666666667
$ws.Range("O16").Value = 0.06
$ws.Range("S16").Value = 0.16
$ws.Range("F17").Value = 0.03061224489795918
$ws.Range("H17").Value = 0.2040816326530612
$ws.Range("I17").Value = 0.09183673469387756
$ws.Range("J17").Value = 0.3801020408163265
$ws.Range("K17").Value = 0.1096938775510204
$ws.Range("M17").Value = 0.01275510204081633
$ws.Range("N17").Value = 0.002551020408163265
$ws.Range("O17").Value = 0.07653061224489796
$ws.Range("S17").Value = 0.09183673469387756
$ws.Range("F18").Value = 0.006329113924050633
$ws.Range("H18").Value = 0.2594936708860759
$ws.Range("I18").Value = 0.05696202531645569
$ws.Range("J18").Value = 0.4367088607594937
$ws.Range("K18").Value = 0.08227848101265822
$ws.Range("M18").Value = 0.0189873417721519
$ws.Range("O18").Value = 0.06962025316455696
$ws.Range("S18").Value = 0.06962025316455696
$ws.Range("F19").Value = 0.01686340640809443
$ws.Range("H19").Value = 0.1913996627318718
$ws.Range("I19").Value = 0.1053962900505902
$ws.Range("J19").Value = 0.3507588532883643
$ws.Range("K19").Value = 0.1247892074198988
$ws.Range("M19").Value = 0.02529510961214165
$ws.Range("N19").Value = 0.002529510961214165
$ws.Range("O19").Value = 0.06661045531197302
$ws.Range("S19").Value = 0.1163575042158516
